# Poster.pptx – split the "Data Analysis" heading run into "Data " + "Analysis"
# (same visible text/formatting, now expressed as two runs) inside the
# "Data Analysis" section-header shape on slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$grp = $s.Shapes.Item("Group 87")
$shp = $grp.GroupItems.Item("Rectangle 89")

$tr = $shp.TextFrame.TextRange

# Re-assigning the text of the leading "Data " substring forces PowerPoint to
# split the single "Data Analysis" run into two runs ("Data " and "Analysis")
# while preserving the existing character formatting (size, bold, fill, fonts).
$tr.Characters(1, 5).Text = "Data "
